# daily auto push: 2025-10-03 22:27 UTC
# This sheet is appended to once per day by an automated job. Add the
# new day's record (2025/10/04, Sat, time=4, ranking=4) as a new row
# right after the current last row of data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
$newRow = $lastRow + 1

# Column A is a date-looking string ("2025/10/04") but, like every
# other row already in the sheet, it must be stored as literal text
# rather than an auto-converted date serial number. A leading
# apostrophe forces Excel to keep it as text.
$ws.Cells.Item($newRow, 1).Value = "'2025/10/04"
$ws.Cells.Item($newRow, 2).Value = "土"
$ws.Cells.Item($newRow, 3).Value = 4
$ws.Cells.Item($newRow, 4).Value = 4
